$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("About")

# New narrow column A (content starts appearing in col A further down already;
# this just gives it an explicit width like the rest of the workbook's columns).
$ws1.Columns.Item(1).ColumnWidth = 7.75

# New note at the bottom of the page, in red, explaining a modelling choice.
$ws1.Range("B24").Value = "use US data"
$ws1.Range("B24").Font.Color = 255
$ws1.Range("B24").Font.Bold = $false

# Scroll / select state left by the editing session.
$ws1.Range("B24").Select()

# ---------------------------------------------------------------------------
# Sheet "EoBSDwEC"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("EoBSDwEC")

# The header over column A used to be a long wrapped label; shorten it to
# "Fuel" and drop the now-unneeded word-wrap / extra row height, matching the
# formatting already used by the other bold headers in this workbook.
$ws2.Range("A1").Value = "Fuel"
$ws1.Range("B7").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Rows.Item(1).RowHeight = 14.45

# Drop the last four fuel rows (kerosene, heavy/residual fuel oil,
# LPG propane/butane, hydrogen) - no longer modelled.
$ws2.Range("A8:D11").EntireRow.Delete()

# Minor column width touch-ups left over from the edit.
$ws2.Columns.Item(1).ColumnWidth = 20.75
